$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows (A, D, E) for worksheet rows 2..57 (A-values 0..55)
$data = @(
    @(0, 0, 0.131),
    @(1, 3, 0.002),
    @(2, 4, 0.012),
    @(3, 5, 0.019),
    @(4, 6, 0.034),
    @(5, 7, 0.045),
    @(6, 8, 0.049),
    @(7, 9, 0.044),
    @(8, 10, 0.045),
    @(9, 11, 0.018),
    @(10, 12, 0.027),
    @(11, 13, 0.033),
    @(12, 14, 0.031),
    @(13, 15, 0.026),
    @(14, 16, 0.034),
    @(15, 17, 0.032),
    @(16, 18, 0.049),
    @(17, 19, 0.035),
    @(18, 20, 0.031),
    @(19, 21, 0.016),
    @(20, 22, 0.018),
    @(21, 23, 0.02),
    @(22, 24, 0.023),
    @(23, 25, 0.018),
    @(24, 26, 0.03),
    @(25, 27, 0.017),
    @(26, 28, 0.014),
    @(27, 29, 0.007),
    @(28, 30, 0.012),
    @(29, 31, 0.011),
    @(30, 32, 0.017),
    @(31, 33, 0.006),
    @(32, 34, 0.017),
    @(33, 35, 0.008),
    @(34, 36, 0.006),
    @(35, 37, 0.006),
    @(36, 38, 0.004),
    @(37, 39, 0.004),
    @(38, 40, 0.003),
    @(39, 41, 0.012),
    @(40, 42, 0.002),
    @(41, 43, 0.006),
    @(42, 44, 0.005),
    @(43, 45, 0.004),
    @(44, 46, 0.002),
    @(45, 47, 0.001),
    @(46, 48, 0.002),
    @(47, 49, 0.002),
    @(48, 50, 0.001),
    @(49, 52, 0.001),
    @(50, 53, 0.002),
    @(51, 54, 0.001),
    @(52, 55, 0.001),
    @(53, 59, 0.001),
    @(54, 60, 0.001),
    @(55, 64, 0.001)
)

$newLamda1 = 33.94444444444444
$newLamda2 = 1.95

# Row 57 is a brand new row; give column A the same style (bold, centered, bordered)
# as the rest of column A's data cells (copy format from A56) before writing its value.
$ws.Range("A56").Copy()
$ws.Range("A57").PasteSpecial(-4122)

$startRow = 2
foreach ($item in $data) {
    $rowIndex = $startRow + [int]$item[0]
    $ws.Cells.Item($rowIndex, 1).Value = $item[0]
    $ws.Cells.Item($rowIndex, 2).Value = $newLamda1
    $ws.Cells.Item($rowIndex, 3).Value = $newLamda2
    $ws.Cells.Item($rowIndex, 4).Value = $item[1]
    $ws.Cells.Item($rowIndex, 5).Value = $item[2]
}
